# Update the Willow TV m3u8 link (cell B2) with a refreshed md5/expires URL.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "https://off1.dunyapurkaraja.com:1686/hls/willowusa.m3u8?md5=1RVIaQIwHsx4YUBN1YJmSw&expires=1742190994"

# Reflect the cell the user last selected/edited.
$ws.Range("B2").Select()

$wb.Save()
